$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("D11").Value = "[1, 0, 1, 0, 1, 0, 0]"
$ws.Range("E11").Value = "['Normal', 'HardwareFault', 'RegulationViolation']"

# Row 26
$ws.Range("D26").Value = "[0, 0, 1, 0, 0, 0, 1]"
$ws.Range("E26").Value = "['HardwareFault', 'SoftwareFault']"

# Row 27
$ws.Range("D27").Value = "[0, 0, 1, 0, 0, 0, 1]"
$ws.Range("E27").Value = "['HardwareFault', 'SoftwareFault']"

# Row 38
$ws.Range("D38").Value = "[0, 0, 1, 0, 0, 0, 0]"
$ws.Range("E38").Value = "['HardwareFault']"

# Row 73
$ws.Range("D73").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal', 'ParamViolation']"

# Row 82
$ws.Range("D82").Value = "[1, 1, 1, 0, 0, 0, 0]"
$ws.Range("E82").Value = "['Normal', 'SurroundingEnvironment', 'HardwareFault']"

$wb.Save()
